$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44188
$ws.Range("M2").Value = 30
$ws.Range("N2").Value = 15000
$ws.Range("O2").Value = 15000
$ws.Range("P2").Value = 15000
$ws.Range("S2").Value = 3000

# Row 3
$ws.Range("D3").Value = 44196
$ws.Range("M3").Value = 56

# Row 4
$ws.Range("D4").Value = 44189

# Row 5
$ws.Range("D5").Value = 44175
$ws.Range("M5").Value = 25

# Row 7
$ws.Range("D7").Value = 44193

# Row 8
$ws.Range("D8").Value = 44186
$ws.Range("M8").Value = 40

# Row 9
$ws.Range("D9").Value = 44179
$ws.Range("M9").Value = 45
$ws.Range("N9").Value = 20000
$ws.Range("O9").Value = 20000
$ws.Range("P9").Value = 20000
$ws.Range("S9").Value = 4000
